$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1837.62
$ws.Range("I15").Value = 1837.62
$ws.Range("K15").Value = 5512.86
$ws.Range("M15").Value = -5343.86
$ws.Range("H33").Value = 320.9
$ws.Range("I33").Value = 172.85715
$ws.Range("K33").Value = 172.85715
$ws.Range("M33").Value = 56.14285000000001
$ws.Range("H51").Value = 10386.4375
$ws.Range("J51").Value = 10279.533
$ws.Range("L51").Value = 10279.533
$ws.Range("N51").Value = -11247.533
$ws.Range("H86").Value = 3612.9092
$ws.Range("I86").Value = 1491.6666
$ws.Range("K86").Value = 1491.6666
$ws.Range("M86").Value = -368.6666
$ws.Range("H89").Value = 3612.9092
$ws.Range("I89").Value = 1491.6666
$ws.Range("K89").Value = 7458.333000000001
$ws.Range("M89").Value = -1842.333000000001
$ws.Range("H111").Value = 7119.4
$ws.Range("I111").Value = 7119.4
$ws.Range("K111").Value = 21358.2
$ws.Range("M111").Value = -18291.2
$ws.Range("H116").Value = 4415.5557
$ws.Range("I116").Value = 3762
$ws.Range("K116").Value = 3762
$ws.Range("M116").Value = -320
$ws.Range("H137").Value = 45382.848
$ws.Range("I137").Value = 86439.69500000001
$ws.Range("J137").Value = 4326
$ws.Range("K137").Value = 259319.085
$ws.Range("L137").Value = 12978
$ws.Range("M137").Value = -256769.085
$ws.Range("N137").Value = -18078
$ws.Range("H138").Value = 3356.9443
$ws.Range("J138").Value = 3323.5625
$ws.Range("L138").Value = 9970.6875
$ws.Range("N138").Value = -20250.6875
$ws.Range("H141").Value = 2268.6
$ws.Range("I141").Value = 2268.6
$ws.Range("K141").Value = 6805.799999999999
$ws.Range("M141").Value = -1625.799999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5337.1904
$ws.Range("I32").Value = 2323.889
$ws.Range("J32").Value = 23417
$ws.Range("K32").Value = 2323.889
$ws.Range("L32").Value = 23417
$ws.Range("M32").Value = -2036.889
$ws.Range("N32").Value = -23991
$ws.Range("H36").Value = 14595.5
$ws.Range("I36").Value = 6969
$ws.Range("J36").Value = 22222
$ws.Range("K36").Value = 6969
$ws.Range("L36").Value = 22222
$ws.Range("M36").Value = -6623
$ws.Range("N36").Value = -22914
$ws.Range("H42").Value = 97515
$ws.Range("J42").Value = 97515
$ws.Range("L42").Value = 97515
$ws.Range("N42").Value = -98487
$ws.Range("H45").Value = 8755.477000000001
$ws.Range("I45").Value = 10116.5625
$ws.Range("K45").Value = 10116.5625
$ws.Range("M45").Value = -9739.5625
$ws.Range("H46").Value = 5152
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5152
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5152
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -5790
$ws.Range("H61").Value = 5011.05
$ws.Range("I61").Value = 4985.316
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 4985.316
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -4773.316
$ws.Range("N61").Value = -5924
$ws.Range("H74").Value = 80869.72
$ws.Range("I74").Value = 80869.72
$ws.Range("K74").Value = 80869.72
$ws.Range("M74").Value = -79995.72
$ws.Range("H77").Value = 80869.72
$ws.Range("I77").Value = 80869.72
$ws.Range("K77").Value = 404348.6
$ws.Range("M77").Value = -399980.6
$ws.Range("H110").Value = 6775.5654
$ws.Range("I110").Value = 6474.65
$ws.Range("J110").Value = 8781.666999999999
$ws.Range("K110").Value = 6474.65
$ws.Range("L110").Value = 8781.666999999999
$ws.Range("M110").Value = -4429.65
$ws.Range("N110").Value = -12871.667
$ws.Range("H124").Value = 29714.5
$ws.Range("J124").Value = 29714.5
$ws.Range("L124").Value = 29714.5
$ws.Range("N124").Value = -39534.5
$ws.Range("H132").Value = 4588.6665
$ws.Range("I132").Value = 3282.6667
$ws.Range("J132").Value = 5241.6665
$ws.Range("K132").Value = 9848.000100000001
$ws.Range("L132").Value = 15724.9995
$ws.Range("M132").Value = -7318.000100000001
$ws.Range("N132").Value = -20784.9995
$ws.Range("H136").Value = 5011.05
$ws.Range("I136").Value = 4985.316
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 14955.948
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -12405.948
$ws.Range("N136").Value = -21600
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 52000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 52000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 52000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -53872
$ws.Range("H77").Value = 52000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 52000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 156000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -165360
$ws.Range("H130").Value = 98000
$ws.Range("J130").Value = 98000
$ws.Range("L130").Value = 98000
$ws.Range("N130").Value = -108040
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 240662.36
$ws.Range("I31").Value = 272544.12
$ws.Range("J31").Value = 4737.2
$ws.Range("K31").Value = 272544.12
$ws.Range("L31").Value = 4737.2
$ws.Range("M31").Value = -272249.12
$ws.Range("N31").Value = -5327.2
$ws.Range("H34").Value = 240662.36
$ws.Range("I34").Value = 272544.12
$ws.Range("J34").Value = 4737.2
$ws.Range("K34").Value = 272544.12
$ws.Range("L34").Value = 4737.2
$ws.Range("M34").Value = -272342.12
$ws.Range("N34").Value = -5141.2
$ws.Range("H58").Value = 4799.9
$ws.Range("I58").Value = 3625.125
$ws.Range("J58").Value = 9499
$ws.Range("K58").Value = 3625.125
$ws.Range("L58").Value = 9499
$ws.Range("M58").Value = -3422.125
$ws.Range("N58").Value = -9905
$ws.Range("H92").Value = 44000
$ws.Range("J92").Value = 44000
$ws.Range("L92").Value = 44000
$ws.Range("N92").Value = -48992
$ws.Range("H99").Value = 487743.38
$ws.Range("I99").Value = 1005861.8
$ws.Range("J99").Value = 16726.637
$ws.Range("K99").Value = 1005861.8
$ws.Range("L99").Value = 16726.637
$ws.Range("M99").Value = -1004363.8
$ws.Range("N99").Value = -19722.637
$ws.Range("H115").Value = 68645
$ws.Range("J115").Value = 68645
$ws.Range("L115").Value = 68645
$ws.Range("N115").Value = -70995
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 487743.38
$ws.Range("I126").Value = 1005861.8
$ws.Range("J126").Value = 16726.637
$ws.Range("K126").Value = 3017585.4
$ws.Range("L126").Value = 50179.91099999999
$ws.Range("M126").Value = -3015115.4
$ws.Range("N126").Value = -55119.91099999999
$ws.Range("H132").Value = 4029.8696
$ws.Range("I132").Value = 4219.4
$ws.Range("K132").Value = 12658.2
$ws.Range("M132").Value = -10128.2
$ws.Range("H136").Value = 4799.9
$ws.Range("I136").Value = 3625.125
$ws.Range("J136").Value = 9499
$ws.Range("K136").Value = 10875.375
$ws.Range("L136").Value = 28497
$ws.Range("M136").Value = -8325.375
$ws.Range("N136").Value = -33597
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6413870
$ws.Range("J68").Value = 4720.8887
$ws.Range("L68").Value = 14162.6661
$ws.Range("N68").Value = -15784.6661
$ws.Range("H71").Value = 6413870
$ws.Range("J71").Value = 4720.8887
$ws.Range("L71").Value = 42487.99830000001
$ws.Range("N71").Value = -50599.99830000001
$ws.Range("H131").Value = 38463692
$ws.Range("J131").Value = 4990
$ws.Range("L131").Value = 14970
$ws.Range("N131").Value = -25050
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 4
$ws.Range("K5").Value = 4
$ws.Range("M5").Value = 108
$ws.Range("H99").Value = 32474.666
$ws.Range("I99").Value = 23624
$ws.Range("J99").Value = 39555.2
$ws.Range("K99").Value = 23624
$ws.Range("L99").Value = 39555.2
$ws.Range("M99").Value = -21378
$ws.Range("N99").Value = -44047.2
$ws.Range("H104").Value = 39633.168
$ws.Range("J104").Value = 39633.168
$ws.Range("L104").Value = 39633.168
$ws.Range("N104").Value = -46621.168
$ws.Range("H107").Value = 51432.3
$ws.Range("I107").Value = 78416.62
$ws.Range("J107").Value = 1318.5714
$ws.Range("K107").Value = 78416.62
$ws.Range("L107").Value = 1318.5714
$ws.Range("M107").Value = -76496.62
$ws.Range("N107").Value = -5158.5714
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1214
$ws.Range("H101").Value = 11992.4
$ws.Range("J101").Value = 11992.4
$ws.Range("L101").Value = 11992.4
$ws.Range("N101").Value = -18482.4
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 54000
$ws.Range("J127").Value = 54000
$ws.Range("L127").Value = 54000
$ws.Range("N127").Value = -63920
$ws.Range("H136").Value = 2299.1365
$ws.Range("I136").Value = 2299.1365
$ws.Range("K136").Value = 6897.4095
$ws.Range("M136").Value = -4347.4095
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 50581.566
$ws.Range("I126").Value = 2009.125
$ws.Range("K126").Value = 6027.375
$ws.Range("M126").Value = -3557.375
$ws.Range("H136").Value = 1429614.4
$ws.Range("J136").Value = 1500
$ws.Range("L136").Value = 4500
$ws.Range("N136").Value = -9600
